$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "7.05", "1.00",
# "0.0000230") retain their exact textual representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.632.23"
$ws.Range("E2").Value = "  -1.45%  "

# Row 3
$ws.Range("D3").Value = "3.025.09"
$ws.Range("E3").Value = "  -1.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "584.13"
$ws.Range("E5").Value = "  -0.74%  "

# Row 6
$ws.Range("D6").Value = "147.33"
$ws.Range("E6").Value = "  -5.20%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "3.023.53"
$ws.Range("E8").Value = "  -1.60%  "

# Row 9
$ws.Range("E9").Value = "  -2.73%  "

# Row 10
$ws.Range("E10").Value = "  -3.88%  "

# Row 11
$ws.Range("E11").Value = "  -1.07%  "

# Row 12
$ws.Range("E12").Value = "  -1.44%  "

# Row 13
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  -2.76%  "

# Row 14
$ws.Range("D14").Value = "34.79"
$ws.Range("E14").Value = "  -5.56%  "

# Row 15
$ws.Range("E15").Value = "  +2.27%  "

# Row 16
$ws.Range("D16").Value = "3.525.48"
$ws.Range("E16").Value = "  -1.59%  "

# Row 17
$ws.Range("D17").Value = "7.05"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18
$ws.Range("D18").Value = "62.618.90"
$ws.Range("E18").Value = "  -1.52%  "

# Row 19
$ws.Range("D19").Value = "3.024.32"
$ws.Range("E19").Value = "  -1.36%  "

# Row 20
$ws.Range("D20").Value = "465.01"
$ws.Range("E20").Value = "  -1.15%  "

# Row 21
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -2.28%  "

# Row 22
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  -2.27%  "

# Row 23
$ws.Range("E23").Value = "  -1.30%  "

# Row 24
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  -3.00%  "

# Row 25
$ws.Range("D25").Value = "80.18"
$ws.Range("E25").Value = "  -0.48%  "

# Row 26
$ws.Range("D26").Value = "12.43"
$ws.Range("E26").Value = "  -2.75%  "

# Row 27
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  -3.04%  "

# Row 28
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("E29").Value = "  +0.39%  "

# Row 30
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  -3.46%  "

# Row 32
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -1.12%  "

# Row 33
$ws.Range("D33").Value = "27.53"
$ws.Range("E33").Value = "  +1.56%  "

# Row 34
$ws.Range("E34").Value = "  -3.81%  "

# Row 35
$ws.Range("E35").Value = "  -0.73%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0795"
$ws.Range("E36").Value = "  -2.91%  "

# Row 37
$ws.Range("E37").Value = "  -3.51%  "

# Row 38
$ws.Range("D38").Value = "2.13"
$ws.Range("E38").Value = "  -3.00%  "

# Row 39
$ws.Range("D39").Value = "50.54"
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("E40").Value = "  -1.63%  "

# Row 41
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").Value = "  -9.87%  "

# Row 42
$ws.Range("D42").Value = "421.00"
$ws.Range("E42").Value = "  -3.46%  "

# Row 43
$ws.Range("E43").Value = "  +1.05%  "

# Row 44
$ws.Range("E44").Value = "  -2.71%  "

# Row 45
$ws.Range("D45").Value = "2.787.88"
$ws.Range("E45").Value = "  -0.22%  "

# Row 46
$ws.Range("E46").Value = "  -1.34%  "

# Row 47
$ws.Range("D47").Value = "37.89"
$ws.Range("E47").Value = "  -6.31%  "

# Row 48
$ws.Range("D48").Value = "129.75"
$ws.Range("E48").Value = "  +0.57%  "

# Row 49
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "24.10"
$ws.Range("E51").Value = "  -3.70%  "
